$wb = $excel.ActiveWorkbook

$colA = @(-1.5,-1.5,-1.5,-1.5,-1.5,-0.4999999999999999,-0.4999999999999999,-0.4999999999999999,-0.4999999999999999,-0.4999999999999999,0.4999999999999999,0.4999999999999999,0.4999999999999999,0.4999999999999999,0.4999999999999999,1.5,1.5,1.5,1.5,1.5,-2,-2,-2,-2,-0.9999999999999999,-0.9999999999999999,-0.9999999999999999,-0.9999999999999999,0,0,0,0,0.9999999999999999,0.9999999999999999,0.9999999999999999,0.9999999999999999,2,2,2,2)
$colB = @(-2,-0.9999999999999999,0,0.9999999999999999,2,-2,-0.9999999999999999,0,0.9999999999999999,2,-2,-0.9999999999999999,0,0.9999999999999999,2,-2,-0.9999999999999999,0,0.9999999999999999,2,-1.5,-0.4999999999999999,0.4999999999999999,1.5,-1.5,-0.4999999999999999,0.4999999999999999,1.5,-1.5,-0.4999999999999999,0.4999999999999999,1.5,-1.5,-0.4999999999999999,0.4999999999999999,1.5,-1.5,-0.4999999999999999,0.4999999999999999,1.5)

$colC1 = @(0.08373889103000398,0.6620963224823134,-1.346262835764032,-0.4943307722557232,0.9080080488278707,-2.416132527868653,0.1126821556674253,-0.01688405974934003,2.00407591647077,3.957524836433219,-3.392024823514296,-0.6376064123218961,2.018028306156258,3.341563561907734,4.153772428731995,-0.9957318400552829,0.2256857690490523,2.098558075797896,1.570175915396323,0.9876082356245576,0.2611324816411806,0.3166552392844378,-2.439768424196422,-1.432261559803352,0.1276654127164068,-0.5594647434738682,0.01515476651823737,0.2644535216577009,1.31840352493853,1.142051091311242,-1.184977370393103,-0.6785372135023864,2.221698517427548,2.645203695189227,1.052308175446992,-0.1249835251422943,1.405454633383981,4.509260994895152,2.618176604340441,0.3236281991753256)
$colC2 = @(0.4235547354964684,1.318300792562311,-1.809008216751274,-0.6756976128341117,0.4293309962818472,-3.011348178320314,0.7479605148049269,-0.6227429467936418,2.149165587918668,4.583632475615934,-4.157983975142266,-1.014258387021291,2.411321557537883,3.654883113363653,4.661222487885713,-0.5373574978510173,0.3854690302727564,2.412356925742293,1.213670282483738,0.5256638275630352,0.1582938849819862,1.27561784729653,-2.977633058578303,-1.736139211348967,0.05049110947680286,-0.8370162697532842,0.2254440975848515,0.5010502203268219,1.67529622925371,1.111493840400597,-1.114108358755421,-1.532992958845985,1.930395494151146,2.962544235498029,1.137673212283299,-0.03358719578136007,1.752126260116206,5.380632818524178,2.204803557752363,0.2035780890178337)
$colC3 = @(-0.07399303219913167,0.533738921045057,-1.149771593746232,-0.2454106212948501,0.8913813849967518,-2.200853693315965,0.1112100541348846,-0.006748851073835924,1.955948625278394,3.63042104778016,-3.011818701263,-0.6067458233117237,1.876152933750652,3.156423019986359,3.798740470875051,-0.8924277848741989,0.2375216025782737,1.916661104140915,1.514929074583702,1.113929239020437,0.1047278650942099,0.2406246687831654,-2.270803639956435,-1.283076272225933,0.0763207633377623,-0.5918768672640989,0.02152886377552686,0.2900251271083729,1.208074491347215,0.9869878184853812,-0.8977256262256325,-0.5018215666187327,2.035394867732024,2.514706662466818,1.134561863745229,-0.06405839566386229,1.447754953416358,4.217343663015215,2.394786495363108,0.5013215664212509)
$colC4 = @(-0.3957668428016821,-0.09439016976558437,0.1972889310505366,0.6389316348028122,1.043348494797261,-0.5248368942899506,-0.04312791381615461,0.4862988709847612,1.079584099309245,1.49293571654793,-0.5048084764601409,0.05952393323149239,0.7688842134362875,1.346861451155469,1.663021951636361,-0.2453654875561058,0.2480898075572119,0.844679586527542,1.249826735534774,1.470066377341473,-0.2493448689506815,-0.2652757297807843,-0.3286289528030609,-0.1537960389324816,0.0853387031419949,-0.01915263475379844,0.1263260740875963,0.3135024829190206,0.5039680614751842,0.5774784758109331,0.288142049356817,0.576692792115349,0.973426859118466,1.114024286770213,1.006634102054557,0.6803517094502072,1.086211937352312,1.505173852818356,1.487546997918387,1.26716868847541)
$colC5 = @(-0.3998262682133235,-0.01459010841871072,-0.3715130475265503,0.2859761851846628,1.046778761770119,-1.333364190470099,-0.1689781646748676,0.3678129985704514,1.594100387995637,2.487782710760611,-1.647062024775973,-0.2919113196789278,1.285858924517725,2.295935018544906,2.711594695330288,-0.700275164980847,0.1826449441317521,1.327798582151998,1.539372257901958,1.499678448565177,-0.1869225499237601,-0.2881258643973917,-1.229349725654885,-0.7447207214452269,-0.01858483360307235,-0.3421310770421605,-0.02356507498586263,0.2704851448398328,0.7664367507608952,0.748381482627263,-0.2873454910569289,0.2470144283454108,1.574837620466897,1.831465566196827,1.144687742871116,0.3086668337075055,1.341853252000011,2.768254863793687,2.098853996599768,1.06093390141675)
$colC6 = @(0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493,0.5173638935743493)

$cArrays = @($colC1, $colC2, $colC3, $colC4, $colC5, $colC6)

for ($s = 1; $s -le 6; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $colC = $cArrays[$s - 1]
    for ($i = 0; $i -lt $colA.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $colA[$i]
        $ws.Cells.Item($row, 2).Value = $colB[$i]
        $ws.Cells.Item($row, 3).Value = $colC[$i]
    }
}

Write-Output "done"
